$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The report lists items alphabetically. A new item ("بامبرز رقم 1" /
# Pampers size 1, price 330.00) needs to be inserted in sorted order between
# "امواس لورد" (row 134) and "جنتيانا " (row 135). That means:
#   - a brand-new row is appended at the bottom of the data block (just above
#     the totals row) to make room,
#   - every existing data row from 135..151 shifts its *values* down by one
#     row (136..152), and
#   - the new product's data is written into row 135.
# The totals row and footer row then move down one row as well, and the
# grand total is bumped by the new item's price.
# ---------------------------------------------------------------------------

# 1) Make room: insert a blank row just above the totals row (152), which
#    pushes the totals row to 153 and the footer row to 154.
$ws.Rows.Item(152).Insert()

# Re-apply the row heights that belong to each row position (Insert() does
# not carry these along automatically).
$ws.Rows.Item(152).RowHeight = 25.5
$ws.Rows.Item(153).RowHeight = 24.75
$ws.Rows.Item(154).RowHeight = 16.5

# Re-create the per-row cell merges for the new data row 152 (matching every
# other data row's merge pattern: A:B, C:G, H:K, L:M, N:O).
$ws.Range("A152:B152").Merge()
$ws.Range("C152:G152").Merge()
$ws.Range("H152:K152").Merge()
$ws.Range("L152:M152").Merge()
$ws.Range("N152:O152").Merge()

# 2) Shift the data values for rows 151 down through 135 into the row below,
#    freeing up row 135 for the new product (working bottom-up so we never
#    overwrite a row before reading it). Column A is just the row's running
#    index (1, 2, 3, ...), not product data, so it is left untouched -- it
#    already reads correctly for every row except the brand-new row 152.
for ($r = 151; $r -ge 135; $r--) {
    $dest = $r + 1

    $cTxt = $ws.Range("C$r").Text
    $hTxt = $ws.Range("H$r").Text
    $lTxt = $ws.Range("L$r").Text
    $nTxt = $ws.Range("N$r").Text
    $pTxt = $ws.Range("P$r").Text
    $qTxt = $ws.Range("Q$r").Text

    $ws.Range("C$dest").Value = "'" + $cTxt
    $ws.Range("H$dest").Value = "'" + $hTxt
    $ws.Range("L$dest").Value = "'" + $lTxt
    $ws.Range("N$dest").Value = "'" + $nTxt
    $ws.Range("P$dest").Value = "'" + $pTxt
    $ws.Range("Q$dest").Value = "'" + $qTxt
}

# Row 152 is brand new, so its running index needs to be filled in explicitly
# (continues the sequence: row 151 -> 145, row 152 -> 146).
$ws.Range("A152").Value = 146

# 3) Write the new product into the now-free row 135 (its running index, 129,
#    was already correct and untouched above).
$ws.Range("C135").Value = "'بامبرز رقم 1"
$ws.Range("H135").Value = "'0:0"
$ws.Range("L135").Value = "'0"
$ws.Range("N135").Value = "'330.00"
$ws.Range("P135").Value = "'330.0000"
$ws.Range("Q135").Value = "'1:0"

# 4) Update the grand total (old total 11624.13 + new item's 330.00).
$oldTotal = $ws.Range("P153").Value2
$ws.Range("P153").Value = $oldTotal + 330

# 5) Refresh the generated-on timestamp in the footer.
$ws.Range("A154").Value = "Sunday, 3 August, 2025 9:28 PM"
